$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new "Decision Tree" values (previously held by row 6)
$ws.Range("A2").Value = "Decision Tree"
$ws.Range("B2").Value = 0.6359790952147639
$ws.Range("C2").Value = 0.507792535872637
$ws.Range("D2").Value = 0.7864561489298842
$ws.Range("E2").Value = 0.3779206446625085
$ws.Range("F2").Value = 0.6743104049374413

# Remove rows 3 through 8 (Random Forest, K-Nearest Neighbors, SVM, Decision Tree(old), Naive Bayes, XGBoost)
$ws.Range("A3:F8").EntireRow.Delete() | Out-Null
